$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Requisitos" entries in rows 24/25 so that the
# "LOM3202 - Circuitos Elétricos (Requisito)" text now appears first (row 24)
# and "LOM3206 - Eletrônica (Indicação de Conjunto)" appears second (row 25).
$reqText = "LOM3202 -  Circuitos Elétricos  (Requisito)`n"
$indText = "LOM3206 -  Eletrônica  (Indicação de Conjunto)`n"

$ws.Range("B24").Value = $reqText
$ws.Range("C24").Value = $reqText

$ws.Range("B25").Value = $indText
$ws.Range("C25").Value = $indText
